$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "x"
$ws.Range("B2").Value = "id-section-bloco-menu"
$ws.Range("E2").Value = "Menu superior fixo"

$ws.Range("B3").Value = "id-article-body"
$ws.Range("D3").Value = "x"
$ws.Range("E3").Value = "representa o body, arquivo principal"

$ws.Columns("B").ColumnWidth = 20.86

$ws.Range("B4").Select()
